$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")

# H8, I8: previously numeric (2.5 / 5) -> now quote-prefixed text "2.500" / "5.000"
$ws.Range("H8").Value = "'2.500"
$ws.Range("I8").Value = "'5.000"

# M8, N8: previously numeric (5 / 10) -> now quote-prefixed text "5.000" / "10.000"
$ws.Range("M8").Value = "'5.000"
$ws.Range("N8").Value = "'10.000"

# R8: device name changed from PRN800 to RS800
$ws.Range("R8").Value = "RS800"

# V8, W8, X8, Y8: recalculated loading values
$ws.Range("V8").Value = 0.405
$ws.Range("W8").Value = 0.724
$ws.Range("X8").Value = 0.405
$ws.Range("Y8").Value = 0.724

# Update selection to R8 (matches sheetView selection in diff)
$ws.Range("R8").Select()
